# Apply weekly report regeneration updates:
# - Refresh "Report Generated On" timestamp
# - Zero out billed amount / pricing totals (no-violation / re-run scenario)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Total Billed Amount summary (C8)
$ws.Range("C8").Value = 0

# Wednesday (07/16/2025) detail pricing + total
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0

# Friday (07/18/2025) detail pricing + total
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
